# Update cryptos price/volume figures per the Sep 18 2023 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '27.366.63'
$ws.Range('E2').Value = '  +2.42%  '

# Row 3
$ws.Range('D3').Value = '1.662.40'
$ws.Range('E3').Value = '  +1.36%  '

# Row 4
$ws.Range('E4').Value = '  -0.44%  '

# Row 5
$ws.Range('D5').Value = '''220.22'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.17%  '

# Row 6
$ws.Range('E6').Value = '  +0.85%  '

# Row 7
$ws.Range('E7').Value = '  -0.43%  '

# Row 8
$ws.Range('D8').Value = '''0.255'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.34%  '

# Row 9
$ws.Range('E9').Value = '  +0.47%  '

# Row 10
$ws.Range('E10').Value = '  +4.71%  '

# Row 11
$ws.Range('D11').Value = '''0.0851'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.05%  '

# Row 12
$ws.Range('D12').Value = '1.894.21'
$ws.Range('E12').Value = '  +1.31%  '

# Row 13
$ws.Range('D13').Value = '1.651.94'
$ws.Range('E13').Value = '  +1.15%  '

# Row 14
$ws.Range('D14').Value = '''4.21'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.41%  '

# Row 15
$ws.Range('D15').Value = '''0.534'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.36%  '

# Row 16
$ws.Range('D16').Value = '''67.32'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.13%  '

# Row 17
$ws.Range('D17').Value = '27.347.47'
$ws.Range('E17').Value = '  +2.32%  '

# Row 18
$ws.Range('D18').Value = '0.0₃0737'

# Row 19
$ws.Range('E19').Value = '  +3.67%  '

# Row 20
$ws.Range('E20').Value = '  -0.38%  '

# Row 22
$ws.Range('E22').Value = '  +2.05%  '

# Row 23
$ws.Range('D23').Value = '''2.52'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +6.44%  '

# Row 25
$ws.Range('D25').Value = '''147.12'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.25%  '

# Row 26
$ws.Range('E26').Value = '  -0.54%  '

# Row 27
$ws.Range('D27').Value = '''7.45'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.00%  '

# Row 28
$ws.Range('E28').Value = '  +1.07%  '

# Row 29
$ws.Range('D29').Value = '''16.06'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.78%  '

# Row 30
$ws.Range('D30').Value = '''0.0515'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.47%  '

# Row 31
$ws.Range('E31').Value = '  +0.88%  '

# Row 32
$ws.Range('E32').Value = '  +0.46%  '

# Row 33
$ws.Range('E33').Value = '  +0.44%  '

# Row 34
$ws.Range('E34').Value = '  +2.47%  '

# Row 35
$ws.Range('D35').Value = '1.265.57'
$ws.Range('E35').Value = '  -1.56%  '

# Row 36
$ws.Range('E36').Value = '  +0.51%  '

# Row 37
$ws.Range('E37').Value = '  +0.12%  '

# Row 38
$ws.Range('E38').Value = '  +0.08%  '

# Row 39
$ws.Range('E39').Value = '  +2.41%  '

# Row 40
$ws.Range('E40').Value = '  -0.34%  '

# Row 41
$ws.Range('D41').Value = '''0.814'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.27%  '

# Row 42
$ws.Range('E42').Value = '  +2.55%  '

# Row 43
$ws.Range('D43').Value = '1.806.41'
$ws.Range('E43').Value = '  +1.51%  '

# Row 44
$ws.Range('E44').Value = '  -4.31%  '

# Row 45
$ws.Range('D45').Value = '''61.85'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.43%  '

# Row 46
$ws.Range('D46').Value = '''92.55'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.85%  '

# Row 47
$ws.Range('E47').Value = '  +1.40%  '

# Row 48
$ws.Range('E48').Value = '  +0.77%  '

# Row 49
$ws.Range('D49').Value = '''0.0987'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.27%  '

# Row 50
$ws.Range('D50').Value = '''7.65'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.64%  '

# Row 51
$ws.Range('E51').Value = '  +0.35%  '
